$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 457.17153048178636
$ws.Range("C2").Value = 387.23121106058704
$ws.Range("D2").Value = 454.71612714123523
$ws.Range("E2").Value = 383.75522269592864

# Row 3 values - C3 removed, D3 added
$ws.Range("B3").Value = 461.94533108605435
$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("D3").Value = 455.6528718572132
$ws.Range("E3").Value = 390.10320212892782

# Selection change
$ws.Range("B1:E3").Select() | Out-Null
